# Reorder the "Recorded By" (column G) comma-separated values so that any
# token equal to "System" (case-insensitive, e.g. "System" or "system")
# is moved to the front of the list, preserving the relative order of the
# remaining tokens and of the System-tokens themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Determine the used range of data rows (row 1 is the header "Recorded By").
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value()

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p.ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -gt 0) {
        $newParts = $systemParts + $otherParts
        $newVal = [string]::Join(", ", $newParts)
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
